$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")

# Update the login password/value on the Login sheet
$wsLogin.Range("B2").Value = 123456789

# Make "Login" the active sheet/tab again (it was "UserData"), and move the
# selection on the Login sheet to B3. Activating "Login" naturally drops
# the previously active sheet's tabSelected flag (on "UserData").
$wsLogin.Activate()
$wsLogin.Range("B3").Select()
